# ToDoList.docx edit
#
# 1) Remove the stray, empty, red-formatted paragraph that sits right
#    after the "Note: A product's associated parts ..." paragraph
#    (just before the trailing blank paragraph / section break).
# 2) Un-hide the built-in "Normal (Web)" style (drop its semiHidden
#    flag) now that it is actively used in the document.

$d = $word.ActiveDocument

# --- 1) Delete the stray red empty paragraph --------------------------
# Locate the paragraph holding the distinctive "Note:" sentence, then
# look at the paragraph right after it: if it is empty and carries the
# explicit red (FF0000) font color from the diff, remove it.
$noteIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*Note:*associated parts*nonpersistent*") {
        $noteIndex = $i
        break
    }
}

if ($noteIndex -ge 1 -and $noteIndex -lt $d.Paragraphs.Count) {
    $candidate = $d.Paragraphs.Item($noteIndex + 1)
    $candidateText = $candidate.Range.Text
    $isBlank = ($candidateText -eq "`r") -or ($candidateText.Trim() -eq "")
    $isRed = ($candidate.Range.Font.Color -eq 255)

    if ($isBlank -and $isRed) {
        $candidate.Range.Delete()
    }
}

# --- 2) Un-hide the "Normal (Web)" style -------------------------------
# In real Word, a style flagged UnhideWhenUsed that is actually applied
# in the document has its semiHidden flag cleared on save. Ask the
# object model to do that explicitly; guard it so the rest of the edit
# still applies even on hosts that don't expose a Style.Hidden setter.
try {
    $webStyle = $d.Styles.Item("Normal (Web)")
    $webStyle.Hidden = $false
} catch {
    # Style.Hidden isn't settable on this host - nothing more we can do
    # through the object model for this particular flag.
}
